$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Year of Treatment" column (column B); everything to its
# right shifts one column to the left.
$ws.Columns.Item(2).Delete()

# Append the ".jamais.jamais" suffix to each (now shifted) header cell.
for ($col = 2; $col -le 7; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value2 = $cell.Value2 + ".jamais.jamais"
}
